# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh: updates the "last updated" timestamp
# in A1, re-sorts a handful of country rows (Armenia/Nigeria,
# Georgia/Burkina Faso, Fiyi/Dominica, Santa Sede/Islas Turcas y Caicos,
# Papua Nueva Guinea/Islas Virgenes Britanicas) and refreshes the numeric
# statistics columns (B:H) for the rows whose figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "Datos actualizados..." timestamp banner -------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 10:05"

# --- Rusia (row 6) - refreshed stats ---------------------------------------
$ws.Cells.Item(6, 1).Value = "Rusia"
$ws.Cells.Item(6, 2).Value = 584680
$ws.Cells.Item(6, 3).Value = 7728
$ws.Cells.Item(6, 4).Value = 339711
$ws.Cells.Item(6, 5).Value = 236858
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 109
$ws.Cells.Item(6, 8).Value = 8111

# --- India (row 7) - refreshed stats ----------------------------------------
$ws.Cells.Item(7, 1).Value = "India"
$ws.Cells.Item(7, 2).Value = 412210
$ws.Cells.Item(7, 3).Value = 483
$ws.Cells.Item(7, 4).Value = 228307
$ws.Cells.Item(7, 5).Value = 170618
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 8
$ws.Cells.Item(7, 8).Value = 13285

# --- Singapur (row 34) - refreshed stats ------------------------------------
$ws.Cells.Item(34, 1).Value = "Singapur"
$ws.Cells.Item(34, 2).Value = 42095
$ws.Cells.Item(34, 3).Value = 262
$ws.Cells.Item(34, 4).Value = 34224
$ws.Cells.Item(34, 5).Value = 7845
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 26

# --- Ucrania (row 38) - refreshed stats -------------------------------------
$ws.Cells.Item(38, 1).Value = "Ucrania"
$ws.Cells.Item(38, 2).Value = 36560
$ws.Cells.Item(38, 3).Value = 735
$ws.Cells.Item(38, 4).Value = 16509
$ws.Cells.Item(38, 5).Value = 19049
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 8
$ws.Cells.Item(38, 8).Value = 1002

# --- Polonia (row 39) - refreshed stats -------------------------------------
$ws.Cells.Item(39, 1).Value = "Polonia"
$ws.Cells.Item(39, 2).Value = 31620
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 16683
$ws.Cells.Item(39, 5).Value = 13591
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 1346

# --- Armenia / Nigeria re-sorted (Armenia now above Nigeria), Armenia refreshed
$ws.Cells.Item(52, 1).Value = "Armenia"
$ws.Cells.Item(52, 2).Value = 20268
$ws.Cells.Item(52, 3).Value = 560
$ws.Cells.Item(52, 4).Value = 9002
$ws.Cells.Item(52, 5).Value = 10916
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 18
$ws.Cells.Item(52, 8).Value = 350

$ws.Cells.Item(53, 1).Value = "Nigeria"
$ws.Cells.Item(53, 2).Value = 19808
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 6718
$ws.Cells.Item(53, 5).Value = 12584
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 506

# --- Hungria (row 89) - refreshed stats -------------------------------------
$ws.Cells.Item(89, 1).Value = "Hungria"
$ws.Cells.Item(89, 2).Value = 4094
$ws.Cells.Item(89, 3).Value = 8
$ws.Cells.Item(89, 4).Value = 2589
$ws.Cells.Item(89, 5).Value = 935
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 570

# --- Estonia (row 105) - refreshed stats ------------------------------------
$ws.Cells.Item(105, 1).Value = "Estonia"
$ws.Cells.Item(105, 2).Value = 1981
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 1764
$ws.Cells.Item(105, 5).Value = 148
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 69

# --- Eslovaquia (row 114) - refreshed stats ---------------------------------
$ws.Cells.Item(114, 1).Value = "Eslovaquia"
$ws.Cells.Item(114, 2).Value = 1587
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = 1447
$ws.Cells.Item(114, 5).Value = 112
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 28

# --- Georgia / Burkina Faso re-sorted (Georgia now above Burkina Faso), Georgia refreshed
$ws.Cells.Item(130, 1).Value = "Georgia"
$ws.Cells.Item(130, 2).Value = 906
$ws.Cells.Item(130, 3).Value = 8
$ws.Cells.Item(130, 4).Value = 755
$ws.Cells.Item(130, 5).Value = 137
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 14

$ws.Cells.Item(131, 1).Value = "Burkina Faso"
$ws.Cells.Item(131, 2).Value = 901
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 810
$ws.Cells.Item(131, 5).Value = 38
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 53

# --- Camboya (row 175) - refreshed stats ------------------------------------
$ws.Cells.Item(175, 1).Value = "Camboya"
$ws.Cells.Item(175, 2).Value = 129
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 127
$ws.Cells.Item(175, 5).Value = 2
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

# --- Fiyi / Dominica re-sorted (Fiyi now above Dominica); figures unchanged
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(202, 2).Value = 18
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 18
$ws.Cells.Item(202, 5).Value = 0
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

$ws.Cells.Item(203, 1).Value = "Dominica"
$ws.Cells.Item(203, 2).Value = 18
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 18
$ws.Cells.Item(203, 5).Value = 0
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0

# --- Santa Sede / Islas Turcas y Caicos re-sorted (Santa Sede now above) ----
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 2).Value = 12
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 2).Value = 12
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 1

# --- Papua Nueva Guinea / Islas Virgenes Britanicas re-sorted (Papua Nueva Guinea now above)
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 2).Value = 8
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1
